# Update countries & provincias Spain
# This workbook lists countries (sheet "Pais") sorted descending by total
# cases (column B). The update refreshes several countries' COVID figures.
# Because the underlying data is kept sorted by total cases, a few rows
# change which country they display:
#   - "Georgia" now has more total cases than "Uganda", "Luxemburgo",
#     "Haiti" and "Gabon", so it moves up from row 113 to row 109, pushing
#     those four countries down by one row each (their own data values
#     travel with them).
#   - "Montserrat" and "Islas Malvinas" are tied on total cases; Montserrat
#     now sorts before Islas Malvinas, so they swap rows 215/216.
# The timestamp in A1 is also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 09:25"

# --- Standalone country data refreshes (no re-sort needed) --------------
# Row 28: Ucrania
$ws.Range("B28").Value = 234584
$ws.Range("C28").Value = 4348
$ws.Range("D28").Value = 103401
$ws.Range("E28").Value = 126663
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 90
$ws.Range("H28").Value = 4520

# Row 63: Armenia
$ws.Range("B63").Value = 53083
$ws.Range("C63").Value = 406
$ws.Range("D63").Value = 44932
$ws.Range("E63").Value = 7161
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 990

# Row 77: Hungria
$ws.Range("B77").Value = 32298
$ws.Range("C77").Value = 818
$ws.Range("D77").Value = 8723
$ws.Range("E77").Value = 22722
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 20
$ws.Range("H77").Value = 853

# Row 79: El Salvador
$ws.Range("E79").Value = 4264
$ws.Range("G79").Value = 4
$ws.Range("H79").Value = 869

# --- Georgia overtakes Uganda / Luxemburgo / Haiti / Gabon ---------------
# Row 109 becomes Georgia with its refreshed totals.
$ws.Range("A109").Value = "Georgia"
$ws.Range("B109").Value = 9245
$ws.Range("C109").Value = 549
$ws.Range("D109").Value = 4887
$ws.Range("E109").Value = 4300
$ws.Range("F109").Value = 0
$ws.Range("G109").Value = 4
$ws.Range("H109").Value = 58

# Row 110 becomes Uganda, keeping Uganda's previous figures.
$ws.Range("A110").Value = "Uganda"
$ws.Range("B110").Value = 8965
$ws.Range("C110").Value = 0
$ws.Range("D110").Value = 5078
$ws.Range("E110").Value = 3805
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 82

# Row 111 becomes Luxemburgo, keeping Luxemburgo's previous figures.
$ws.Range("A111").Value = "Luxemburgo"
$ws.Range("B111").Value = 8925
$ws.Range("C111").Value = 0
$ws.Range("D111").Value = 7793
$ws.Range("E111").Value = 1005
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 127

# Row 112 becomes Haiti, keeping Haiti's previous figures.
$ws.Range("A112").Value = "Haiti"
$ws.Range("B112").Value = 8827
$ws.Range("C112").Value = 0
$ws.Range("D112").Value = 6992
$ws.Range("E112").Value = 1606
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 229

# Row 113 becomes Gabon, keeping Gabon's previous figures.
$ws.Range("A113").Value = "Gabon"
$ws.Range("B113").Value = 8808
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 8135
$ws.Range("E113").Value = 619
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 54

# --- Montserrat overtakes Islas Malvinas (tie-break reorder) ------------
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
